# Updates cryptos list price/volume data (and restores the
# Polygon/EnergySwap row ordering) per the target commit.
#
# Every write goes through a scratch cell (Z1) that is force-
# formatted as Text ("@") before the value is assigned, then
# copied across with PasteSpecial(xlPasteValues). This keeps
# every updated cell a plain text value (matching the original
# inline-string cells, e.g. "0.990" keeps its trailing zero,
# "542.82" is not reinterpreted as a number, etc.) while leaving
# the destination cell's own style/number-format untouched --
# a direct "$ws.Range(...).Value = ..." would silently coerce
# numeric-looking strings to numbers and/or stamp a new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")
$xlPasteValues = -4163

$helper.NumberFormat = '@'
$helper.Value = '57.450.11'
$helper.Copy()
$ws.Range('D2').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.01%  '
$helper.Copy()
$ws.Range('E2').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '2.329.54'
$helper.Copy()
$ws.Range('D3').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.14%  '
$helper.Copy()
$ws.Range('E3').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.12%  '
$helper.Copy()
$ws.Range('E4').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '542.82'
$helper.Copy()
$ws.Range('D5').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +6.38%  '
$helper.Copy()
$ws.Range('E5').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '135.15'
$helper.Copy()
$ws.Range('D6').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.31%  '
$helper.Copy()
$ws.Range('E6').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.992'
$helper.Copy()
$ws.Range('D7').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  -0.66%  '
$helper.Copy()
$ws.Range('E7').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.535'
$helper.Copy()
$ws.Range('D8').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.72%  '
$helper.Copy()
$ws.Range('E8').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '2.361.54'
$helper.Copy()
$ws.Range('D9').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.45%  '
$helper.Copy()
$ws.Range('E9').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.68%  '
$helper.Copy()
$ws.Range('E10').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.11%  '
$helper.Copy()
$ws.Range('E11').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '5.38'
$helper.Copy()
$ws.Range('D12').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.80%  '
$helper.Copy()
$ws.Range('E12').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.355'
$helper.Copy()
$ws.Range('D13').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +5.58%  '
$helper.Copy()
$ws.Range('E13').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '2.752.83'
$helper.Copy()
$ws.Range('D14').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.47%  '
$helper.Copy()
$ws.Range('E14').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '23.54'
$helper.Copy()
$ws.Range('D15').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.56%  '
$helper.Copy()
$ws.Range('E15').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '57.512.75'
$helper.Copy()
$ws.Range('D16').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.16%  '
$helper.Copy()
$ws.Range('E16').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.08%  '
$helper.Copy()
$ws.Range('E17').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '2.341.49'
$helper.Copy()
$ws.Range('D18').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.69%  '
$helper.Copy()
$ws.Range('E18').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '337.95'
$helper.Copy()
$ws.Range('D19').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +5.00%  '
$helper.Copy()
$ws.Range('E19').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.65%  '
$helper.Copy()
$ws.Range('E20').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '4.21'
$helper.Copy()
$ws.Range('D21').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.05%  '
$helper.Copy()
$ws.Range('E21').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +3.64%  '
$helper.Copy()
$ws.Range('E22').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.994'
$helper.Copy()
$ws.Range('D23').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  -0.40%  '
$helper.Copy()
$ws.Range('E23').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '61.94'
$helper.Copy()
$ws.Range('D24').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.86%  '
$helper.Copy()
$ws.Range('E24').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +5.33%  '
$helper.Copy()
$ws.Range('E25').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '8.48'
$helper.Copy()
$ws.Range('D26').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  -1.96%  '
$helper.Copy()
$ws.Range('E26').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.995'
$helper.Copy()
$ws.Range('D27').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  -0.41%  '
$helper.Copy()
$ws.Range('E27').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '1.43'
$helper.Copy()
$ws.Range('D28').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +9.67%  '
$helper.Copy()
$ws.Range('E28').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +5.66%  '
$helper.Copy()
$ws.Range('E29').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '170.52'
$helper.Copy()
$ws.Range('D30').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.95%  '
$helper.Copy()
$ws.Range('E30').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.0₃0738'
$helper.Copy()
$ws.Range('D31').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.95%  '
$helper.Copy()
$ws.Range('E31').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.78%  '
$helper.Copy()
$ws.Range('E32').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '18.59'
$helper.Copy()
$ws.Range('D33').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.01%  '
$helper.Copy()
$ws.Range('E33').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +15.64%  '
$helper.Copy()
$ws.Range('E34').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.998'
$helper.Copy()
$ws.Range('D35').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  -0.10%  '
$helper.Copy()
$ws.Range('E35').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.990'
$helper.Copy()
$ws.Range('D36').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  -0.76%  '
$helper.Copy()
$ws.Range('E36').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '4.14'
$helper.Copy()
$ws.Range('D37').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +5.91%  '
$helper.Copy()
$ws.Range('E37').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.33%  '
$helper.Copy()
$ws.Range('E38').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +4.90%  '
$helper.Copy()
$ws.Range('E39').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '39.37'
$helper.Copy()
$ws.Range('D40').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.27%  '
$helper.Copy()
$ws.Range('E40').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '148.71'
$helper.Copy()
$ws.Range('D41').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  -1.45%  '
$helper.Copy()
$ws.Range('E41').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.378'
$helper.Copy()
$ws.Range('D42').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +0.98%  '
$helper.Copy()
$ws.Range('E42').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.27%  '
$helper.Copy()
$ws.Range('E43').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '284.04'
$helper.Copy()
$ws.Range('D44').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.11%  '
$helper.Copy()
$ws.Range('E44').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.18%  '
$helper.Copy()
$ws.Range('E45').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '19.21'
$helper.Copy()
$ws.Range('D46').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +7.16%  '
$helper.Copy()
$ws.Range('E46').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.59%  '
$helper.Copy()
$ws.Range('E47').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.561'
$helper.Copy()
$ws.Range('D48').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.10%  '
$helper.Copy()
$ws.Range('E48').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +2.74%  '
$helper.Copy()
$ws.Range('E49').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = 'EnergySwap'
$helper.Copy()
$ws.Range('B50').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$helper.Copy()
$ws.Range('C50').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '17.57'
$helper.Copy()
$ws.Range('D50').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +3.31%  '
$helper.Copy()
$ws.Range('E50').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = 'Polygon'
$helper.Copy()
$ws.Range('B51').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$helper.Copy()
$ws.Range('C51').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '0.382'
$helper.Copy()
$ws.Range('D51').PasteSpecial($xlPasteValues)
$helper.NumberFormat = '@'
$helper.Value = '  +1.23%  '
$helper.Copy()
$ws.Range('E51').PasteSpecial($xlPasteValues)

$helper.Clear()
$excel.CutCopyMode = 0
